# Remove the "new_fact_currency_rate" row from the pipeline table definition.
# This row (row 30: stage=raw, table_type=raw, database=raw_ad_works_dw,
# table=new_fact_currency_rate) duplicated fact_currency_rate and is removed
# as part of fixing the pipeline file names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the entire row 30. This shifts all subsequent rows up by one, which
# also updates the worksheet dimension and the frozen-pane state.
$ws.Rows.Item(30).Delete() | Out-Null

# The hidden _xlnm._FilterDatabase defined name still points at the old
# (pre-delete) range, so shift it up by one row to match.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$F`$30:`$G`$55"
    }
}

# Leave the active selection on the row that used to hold the deleted entry,
# matching the cursor/selection Excel leaves behind after a row delete.
$ws.Rows.Item(30).Select() | Out-Null
